$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "80.142.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.223.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.82%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +7.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "641.38"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.271"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +31.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +10.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.221.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.622"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +42.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000274"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +43.52%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.46"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.816.22"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.87"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +14.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.862.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.206.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.05"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +28.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.43"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +19.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.37"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +23.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.90"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +13.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.13"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.366.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +12.87%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000127"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +18.77%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.34"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +13.61%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.81"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +15.81%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.54"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +11.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.157"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +28.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.58"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +15.11%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +20.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.417"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +11.12%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +15.05%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "163.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.32"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "193.36"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.85"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +13.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +14.44%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.41%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.804"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +13.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.651"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.32%  "
